$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 37; this shifts the existing rows 37..102
# down to 38..103 (and grows the used range to A1:T103).
$ws.Rows.Item(37).Insert()

# Populate the newly inserted row 37 with the new weekly record
# (Papaya, Vega Modelo de Temuco).
$ws.Cells.Item(37, 1).Value = 10
$ws.Cells.Item(37, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(37, 3).Value = "La Araucanía"
$ws.Cells.Item(37, 4).Value = 45028
$ws.Cells.Item(37, 5).Value = 9
$ws.Cells.Item(37, 6).Value = "Fruta"
$ws.Cells.Item(37, 7).Value = 100108
$ws.Cells.Item(37, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(37, 9).Value = 100108004
$ws.Cells.Item(37, 10).Value = "Papaya"
$ws.Cells.Item(37, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(37, 12).Value = "Primera"
$ws.Cells.Item(37, 13).Value = 200
$ws.Cells.Item(37, 14).Value = 3500
$ws.Cells.Item(37, 15).Value = 3500
$ws.Cells.Item(37, 16).Value = 3500
$ws.Cells.Item(37, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(37, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(37, 19).Value = 3500
$ws.Cells.Item(37, 20).Value = 1
